# Michael Jackson's Moonwalker plan.xlsx -- "Updated set to better align with writing policy"
#
# The commit reworks a batch of Achievement descriptions on the "Achievements"
# sheet (column F) to remove parenthetical asides / bracketed tags in favour
# of plain comma- and "on "-joined phrasing. All the dependent sheets
# (Checklist, Text, Leaderboards, Stats) pull these strings via formulas
# (e.g. =Achievements!F3), so updating the source cells here is sufficient
# for them to recompute.
#
# The workbook's active-tab/selection also moved from "Regions Test"!G6:R6
# to "Achievements"!F34 (with the visible window scrolled to show row ~19).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Achievements")

$replacements = @{
    "F3"  = "Perform a Moonwalk for 5 seconds, counts time while moving, resets when Moonwalk ends"
    "F10" = "Moonwalk through a spider web in stage 4 any round, start Moonwalk before entering web, no turning around"
    "F11" = "Break all nine computers and servers in stage 5-1"
    "F12" = "Moonwalk over the entire bottom floor conveyor in stage 5-2 travelling against the direction of conveyor, left to right without turning around"
    "F15" = "Complete Stage 1-1 under the par time 2m00s on Normal+"
    "F16" = "Complete Stage 1-2 under the par time 2m00s on Normal+"
    "F17" = "Complete Stage 1-3 under the par time 2m00s on Normal+"
    "F18" = "Complete Stage 2-1 under the par time 3m00s on Normal+"
    "F19" = "Complete Stage 2-2 under the par time 3m00s on Normal+"
    "F20" = "Complete Stage 2-3 under the par time 3m00s on Normal+"
    "F21" = "Complete Stage 3-1 under the par time 4m00s on Normal+"
    "F22" = "Complete Stage 3-2 under the par time 4m00s on Normal+"
    "F23" = "Complete Stage 3-3 under the par time 4m00s on Normal+"
    "F24" = "Complete Stage 4-1 under the par time 5m00s on Normal+"
    "F25" = "Complete Stage 4-2 under the par time 5m00s on Normal+"
    "F26" = "Complete Stage 4-3 under the par time 5m00s on Normal+"
    "F27" = "Complete Stage 5-1 under the par time 6m00s on Normal+"
    "F28" = "Complete Stage 5-2 under the par time 6m00s on Normal+"
    "F52" = "Clear Stage 6-1 on Hard"
    "F53" = "Clear Stage 6-1 without using a continue on Normal+"
}

foreach ($addr in $replacements.Keys) {
    $ws.Range($addr).Value = $replacements[$addr]
}

# Move the active tab / selection from Regions Test!G6:R6 to Achievements!F34
$ws.Activate()
$ws.Range("F34").Select()
